# Trade #104 closed at 2026-02-17 15:58:43 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Summary"
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Range("B3").Value = 1200        # Current Capital
$summary.Range("B4").Value = -0.01       # Total P&L $
$summary.Range("B6").Value = 104         # Total Trades
$summary.Range("B7").Value = 40          # Winning Trades
$summary.Range("B9").Value = 38.46       # Win Rate %

# ---------------------------------------------------------------------
# Sheet 2: "Strategy Status" - row 4 is the MarketMaking strategy
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item(2)
$status.Range("C4").Value = 100          # Capital
$status.Range("D4").Value = 104          # Trades
$status.Range("E4").Value = -0.01        # P&L $
$status.Range("F4").Value = 0            # P&L %
$status.Range("G4").Value = 38.46        # Win Rate %

# ---------------------------------------------------------------------
# Sheets 3 & 4: "All Trades" and "MarketMaking" - append trade #104 as
# a new row (row 105). Both sheets carry identical trade logs.
# ---------------------------------------------------------------------
$newRow = @{
    A = 104
    B = "2026-02-17"
    C = "15:58:36"
    D = "MarketMaking"
    E = "DOWN"
    F = 0.96
    G = 0.97
    H = "CLOSED"
    I = 1.0417
    J = 0.01
    K = 100
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.11
}

foreach ($sheetIndex in 3, 4) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    $cellB = $ws.Range("B105")
    $cellB.NumberFormat = "@"
    $cellB.Value = $newRow.B

    $cellC = $ws.Range("C105")
    $cellC.NumberFormat = "@"
    $cellC.Value = $newRow.C

    $ws.Range("A105").Value = $newRow.A
    $ws.Range("D105").Value = $newRow.D
    $ws.Range("E105").Value = $newRow.E
    $ws.Range("F105").Value = $newRow.F
    $ws.Range("G105").Value = $newRow.G
    $ws.Range("H105").Value = $newRow.H
    $ws.Range("I105").Value = $newRow.I
    $ws.Range("J105").Value = $newRow.J
    $ws.Range("K105").Value = $newRow.K
    $ws.Range("L105").Value = $newRow.L
    $ws.Range("M105").Value = $newRow.M
    $ws.Range("N105").Value = $newRow.N
    $ws.Range("O105").Value = $newRow.O
    $ws.Range("P105").Value = $newRow.P
    $ws.Range("Q105").Value = $newRow.Q
}
